$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("E2").Value = 3
$ws.Range("F2").Value = 1
$ws.Range("G2").Value = 2.276052666666667
$ws.Range("H2").Value = 6.828158
$ws.Range("I2").Value = 0.005247614157263819
$ws.Range("J2").Value = 0.005247614157263819
$ws.Range("M2").Value = 8.820647333333334
$ws.Range("N2").Value = 26.461942
$ws.Range("O2").Value = 0.06415146660411865
$ws.Range("P2").Value = 0.06415146660411865
$ws.Range("Q2").Value = 20.07625788475956
$ws.Range("R2").Value = 180.686320962836
$ws.Range("S2").Value = 0.0003366421443610101
$ws.Range("T2").Value = 0.0003366421443610101
$ws.Range("E3").Value = 3
$ws.Range("F3").Value = 1
$ws.Range("G3").Value = 2.276052666666667
$ws.Range("H3").Value = 6.828158
$ws.Range("I3").Value = 0.005247614157263819
$ws.Range("J3").Value = 0.005247614157263819
$ws.Range("O3").Value = 0.3979101621202897
$ws.Range("P3").Value = 0.3979101621202898
$ws.Range("Q3").Value = 124.5263351341767
$ws.Range("R3").Value = 1120.73701620759
$ws.Range("S3").Value = 0.002088079000061574
$ws.Range("T3").Value = 0.002088079000061574
$ws.Range("E4").Value = 3
$ws.Range("F4").Value = 1
$ws.Range("G4").Value = 2.276052666666667
$ws.Range("H4").Value = 6.828158
$ws.Range("I4").Value = 0.005247614157263819
$ws.Range("J4").Value = 0.005247614157263819
$ws.Range("M4").Value = 21.90816333333333
$ws.Range("N4").Value = 65.72449
$ws.Range("O4").Value = 0.1593353362087987
$ws.Range("P4").Value = 0.1593353362087987
$ws.Range("Q4").Value = 49.86413357660223
$ws.Range("R4").Value = 448.7772021894201
$ws.Range("S4").Value = 0.0008361303660416822
$ws.Range("T4").Value = 0.0008361303660416823
$ws.Range("E5").Value = 3
$ws.Range("F5").Value = 1
$ws.Range("G5").Value = 2.276052666666667
$ws.Range("H5").Value = 6.828158
$ws.Range("I5").Value = 0.005247614157263819
$ws.Range("J5").Value = 0.005247614157263819
$ws.Range("M5").Value = 52.056859
$ws.Range("N5").Value = 156.170577
$ws.Range("O5").Value = 0.3786030350667928
$ws.Range("P5").Value = 0.3786030350667929
$ws.Range("Q5").Value = 118.4841527452407
$ws.Range("R5").Value = 1066.357374707166
$ws.Range("S5").Value = 0.001986762646799552
$ws.Range("T5").Value = 0.001986762646799553
$ws.Range("I6").Value = 0.1062533062835484
$ws.Range("J6").Value = 0.1062533062835484
$ws.Range("M6").Value = 8.820647333333334
$ws.Range("N6").Value = 26.461942
$ws.Range("O6").Value = 0.06415146660411865
$ws.Range("P6").Value = 0.06415146660411865
$ws.Range("Q6").Value = 406.5025960615072
$ws.Range("R6").Value = 3658.523364553564
$ws.Range("S6").Value = 0.006816305429626244
$ws.Range("T6").Value = 0.006816305429626245
$ws.Range("I7").Value = 0.1062533062835484
$ws.Range("J7").Value = 0.1062533062835484
$ws.Range("O7").Value = 0.3979101621202897
$ws.Range("P7").Value = 0.3979101621202898
$ws.Range("S7").Value = 0.04227927032910354
$ws.Range("T7").Value = 0.04227927032910355
$ws.Range("I8").Value = 0.1062533062835484
$ws.Range("J8").Value = 0.1062533062835484
$ws.Range("M8").Value = 21.90816333333333
$ws.Range("N8").Value = 65.72449
$ws.Range("O8").Value = 0.1593353362087987
$ws.Range("P8").Value = 0.1593353362087987
$ws.Range("Q8").Value = 1009.645316652065
$ws.Range("R8").Value = 9086.807849868581
$ws.Range("S8").Value = 0.01692990627998564
$ws.Range("T8").Value = 0.01692990627998564
$ws.Range("I9").Value = 0.1062533062835484
$ws.Range("J9").Value = 0.1062533062835484
$ws.Range("M9").Value = 52.056859
$ws.Range("N9").Value = 156.170577
$ws.Range("O9").Value = 0.3786030350667928
$ws.Range("P9").Value = 0.3786030350667929
$ws.Range("Q9").Value = 2399.058428097359
$ws.Range("R9").Value = 21591.52585287623
$ws.Range("S9").Value = 0.04022782424483295
$ws.Range("T9").Value = 0.04022782424483296
$ws.Range("G10").Value = 41.187613
$ws.Range("H10").Value = 123.562839
$ws.Range("I10").Value = 0.09496120377532416
$ws.Range("J10").Value = 0.09496120377532417
$ws.Range("M10").Value = 8.820647333333334
$ws.Range("N10").Value = 26.461942
$ws.Range("O10").Value = 0.06415146660411865
$ws.Range("P10").Value = 0.06415146660411865
$ws.Range("Q10").Value = 363.3014087748153
$ws.Range("R10").Value = 3269.712678973338
$ws.Range("S10").Value = 0.006091900492679613
$ws.Range("T10").Value = 0.006091900492679614
$ws.Range("G11").Value = 41.187613
$ws.Range("H11").Value = 123.562839
$ws.Range("I11").Value = 0.09496120377532416
$ws.Range("J11").Value = 0.09496120377532417
$ws.Range("O11").Value = 0.3979101621202897
$ws.Range("P11").Value = 0.3979101621202898
$ws.Range("Q11").Value = 2253.437530215955
$ws.Range("R11").Value = 20280.93777194359
$ws.Range("S11").Value = 0.03778602798937711
$ws.Range("T11").Value = 0.03778602798937711
$ws.Range("G12").Value = 41.187613
$ws.Range("H12").Value = 123.562839
$ws.Range("I12").Value = 0.09496120377532416
$ws.Range("J12").Value = 0.09496120377532417
$ws.Range("M12").Value = 21.90816333333333
$ws.Range("N12").Value = 65.72449
$ws.Range("O12").Value = 0.1593353362087987
$ws.Range("P12").Value = 0.1593353362087987
$ws.Range("Q12").Value = 902.3449529141234
$ws.Range("R12").Value = 8121.104576227111
$ws.Range("S12").Value = 0.01513067533033352
$ws.Range("T12").Value = 0.01513067533033352
$ws.Range("G13").Value = 41.187613
$ws.Range("H13").Value = 123.562839
$ws.Range("I13").Value = 0.09496120377532416
$ws.Range("J13").Value = 0.09496120377532417
$ws.Range("M13").Value = 52.056859
$ws.Range("N13").Value = 156.170577
$ws.Range("O13").Value = 0.3786030350667928
$ws.Range("P13").Value = 0.3786030350667929
$ws.Range("Q13").Value = 2144.097762487567
$ws.Range("R13").Value = 19296.8798623881
$ws.Range("S13").Value = 0.03595259996293391
$ws.Range("T13").Value = 0.03595259996293392
$ws.Range("G14").Value = 344.1819356666667
$ws.Range("H14").Value = 1032.545807
$ws.Range("I14").Value = 0.7935378757838636
$ws.Range("J14").Value = 0.7935378757838637
$ws.Range("M14").Value = 8.820647333333334
$ws.Range("N14").Value = 26.461942
$ws.Range("O14").Value = 0.06415146660411865
$ws.Range("P14").Value = 0.06415146660411865
$ws.Range("Q14").Value = 3035.907473019688
$ws.Range("R14").Value = 27323.16725717719
$ws.Range("S14").Value = 0.05090661853745178
$ws.Range("T14").Value = 0.05090661853745179
$ws.Range("G15").Value = 344.1819356666667
$ws.Range("H15").Value = 1032.545807
$ws.Range("I15").Value = 0.7935378757838636
$ws.Range("J15").Value = 0.7935378757838637
$ws.Range("O15").Value = 0.3979101621202897
$ws.Range("P15").Value = 0.3979101621202898
$ws.Range("Q15").Value = 18830.72201959458
$ws.Range("R15").Value = 169476.4981763512
$ws.Range("S15").Value = 0.3157567848017475
$ws.Range("T15").Value = 0.3157567848017476
$ws.Range("G16").Value = 344.1819356666667
$ws.Range("H16").Value = 1032.545807
$ws.Range("I16").Value = 0.7935378757838636
$ws.Range("J16").Value = 0.7935378757838637
$ws.Range("M16").Value = 21.90816333333333
$ws.Range("N16").Value = 65.72449
$ws.Range("O16").Value = 0.1593353362087987
$ws.Range("P16").Value = 0.1593353362087987
$ws.Range("Q16").Value = 7540.394062968159
$ws.Range("R16").Value = 67863.54656671343
$ws.Range("S16").Value = 0.1264386242324378
$ws.Range("T16").Value = 0.1264386242324378
$ws.Range("G17").Value = 344.1819356666667
$ws.Range("H17").Value = 1032.545807
$ws.Range("I17").Value = 0.7935378757838636
$ws.Range("J17").Value = 0.7935378757838637
$ws.Range("M17").Value = 52.056859
$ws.Range("N17").Value = 156.170577
$ws.Range("O17").Value = 0.3786030350667928
$ws.Range("R17").Value = 161253.2744581206
$ws.Range("S17").Value = 0.3004358482122264
$ws.Range("T17").Value = 0.3004358482122265
$ws.Range("P17").Value = 0.3786030350667929
